# Update Rspo2-Lgr6.xlsx with newly recomputed TPM-based NATMI statistics.
#
# The underlying analysis was rerun with new TPM values. With the updated
# data, the "ECs" target-cluster row no longer qualifies, and the
# remaining (surviving) row corresponds to the "MuSCs" target cluster with
# freshly recalculated receptor/edge statistics (columns K:T). The old
# second data row (originally the MuSCs row) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that remains (row 2) now represents the MuSCs target cluster.
$ws.Range("D2").Value = "MuSCs"

# Recomputed receptor / edge statistics for the surviving row.
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4426103333333333
$ws.Range("N2").Value = 1.327831
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.07552289625022222
$ws.Range("R2").Value = 0.679706066252
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Drop the old third row (previously the MuSCs row) entirely; the sheet
# now only has the header row plus one data row.
$ws.Rows.Item(3).Delete()
